# projektjournal-tnt.xlsx update
# - Zeitplan / Projektjournal: revise texts for rows 33-34, add new journal
#   entries in rows 35-37, update column width + navigation/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------
# Row 33 - revise "Realisation geplant" entry (Beschreibung column).
# The Diskussion column (F33) keeps its existing text - no edit needed,
# it just shifts down in the shared-string table automatically.
# ---------------------------------------------------------------------
$ws.Range("E33").Value = "Wir haben die Arbeitschritte unter uns aufgeteilt und das Datenbank-Modell umgesetzt."

# ---------------------------------------------------------------------
# Row 34 - "Realisation begonnen" entry gets a longer description and
# now wraps onto two lines. D34 keeps its existing text, no edit needed.
# ---------------------------------------------------------------------
$ws.Range("E34").Value = "Heute haben wir mit dem GUI begonnen. Zudem haben wir mit der Umsetzung der Navigation begonnen."
$ws.Rows("34").RowHeight = 25.5

# ---------------------------------------------------------------------
# Rows 35-37 used to be empty placeholder rows - fill them with the new
# journal entries. Copy the date-cell format from row 34 first so the
# date style (numFmt) matches the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 35
$ws.Range("A35").Value = 43517
$ws.Range("B35").Value = 300
$ws.Range("C35").Value = "Elias, Nico"
$ws.Range("D35").Value = "GUI und Datenbank API Service weiterführen"
$ws.Range("E35").Value = "Nico und Elias haben in der Freizeit weiter am GUI und an den Datenbank API Services weitergearbeitet."
$ws.Range("F35").Value = "Elias hatte zuerst Schwierigkeiten, da dieser sich nicht mit Laravel auskennt. Nico konnte Elias dabei noch weiterhelfen."
$ws.Rows("35").RowHeight = 25.5

# Row 36 (Diskussion column F36 was authored before the Beschreibung
# column E36, so write F36 first to mirror the original shared-string
# insertion order)
$ws.Range("A36").Value = 43522
$ws.Range("B36").Value = 270
$ws.Range("C36").Value = "Elias, Nico, Tim"
$ws.Range("D36").Value = "GUI und Datenbank API Service abschliessen."
$ws.Range("F36").Value = "Durch Nico konnte Elias und Tim noch einiges über PHP Laravel lernen."
$ws.Range("E36").Value = "Heute konnte Nico das GUI und die Datenbank API Service abschliessen. Zudem haben wir noch die Arbeitsaufteilung der Umsetzung der Navigation vollzogen"
$ws.Rows("36").RowHeight = 38.25

# Row 37
$ws.Range("A37").Value = 43509
$ws.Range("B37").Value = 155
$ws.Range("C37").Value = "Elias, Nico, Tim"
$ws.Range("D37").Value = "Administratives nachführen"
$ws.Range("E37").Value = "Heute konnten wir, wegen eines Wegfallens durch Tim und Nico, nicht weiter an der Realisation arbeiten. Elias hat einige administrative Arbeiten, wie den Zeitplan und das Arbeitsjournal, nachführen können."
$ws.Range("F37").Value = "Der Wegfall von Tim und Nico, könnte dafür sorgen, dass dies erst nächste Woche aufarbeiten können. Dies Verschiebt unsere Planung, fällt jedoch nicht sehr gravierend auf."
$ws.Rows("37").RowHeight = 51

# ---------------------------------------------------------------------
# Widen column D slightly to fit the new activity texts
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 40.83

# ---------------------------------------------------------------------
# Move the view / selection like the author left it (scrolled further
# down, cursor parked on D38)
# ---------------------------------------------------------------------
$ws.Range("D38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.Calculate()
